# dy_qT/expdata/90002.xlsx — "fiducial=False, %sys_c -> %norm_c"
#
# 1) Rename the header in column N from "sys,luminosity uncertainty +" to "%norm_c".
# 2) Replace the ymin/ymax (F,G) and etaMin/etaMax (R,S) numeric sentinel values
#    (-1000 / 1000) with boolean FALSE for every data row (2-83), i.e. turn off
#    those fiducial-cut ranges.
# 3) Give column N (14) an explicit width.
# 4) Reset the sheet view back to the top-left / A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 1) Shared-string header rename --------------------------------------
$ws.Range("N1").Value = "%norm_c"

# -- 2) Convert F,G,R,S columns to boolean FALSE for rows 2..83 ----------
for ($r = 2; $r -le 83; $r++) {
    $ws.Cells.Item($r, 6).Value  = $false   # F -> ymin
    $ws.Cells.Item($r, 7).Value  = $false   # G -> ymax
    $ws.Cells.Item($r, 18).Value = $false   # R -> etaMin
    $ws.Cells.Item($r, 19).Value = $false   # S -> etaMax
}

# -- 3) Column N width -----------------------------------------------------
$ws.Columns.Item(14).ColumnWidth = 10.2

# -- 4) Reset the view to A1 / top-left ------------------------------------
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
